$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (B2, C2, D2 change; E2 stays the same)
$ws.Range("B2").Value = "5010754-58.2017.8.21.0001"
$ws.Range("C2").Value = "0196807-38.2017.8.21.0001"
$ws.Range("D2").Value = "CIV.04574.01"

# Copy the formatting of A2 (bold, centered, bordered) onto the new A3:A5 cells
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)

# Add new rows 3-5
$data = @(
    @(1, "5001221-76.2018.8.21.0054", "0066578-85.2019.8.21.9000", "CIV.35589.02", "originario_principal"),
    @(2, "5009614-96.2011.8.21.0001", "0420415-91.2011.8.21.0001", "CIV.30745.01", "originario_principal"),
    @(3, "5033806-25.2013.8.21.0001", "0007850-16.2013.8.21.3001", "CIV.27994.01", "originario_principal")
)

$row = 3
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
